$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Training Dashboard")

# --- Font/style update -------------------------------------------------
# The title (row 1) used to be bold 14pt; it now shares the same bold,
# default-size (11pt) white font that the header row (row 2) also now uses.
$titleCell = $ws.Range("A1")
$titleCell.Font.Size = 11
$titleCell.Font.Color = 16777215

$headerRange = $ws.Range("A2:K2")
$headerRange.Font.Color = 16777215

# --- Data updates (rows 3-5): refreshed "PERIOD TO EXPIRE" (numeric) and
# "LAST UPDATE" (text date label) columns.
# Numbers can be written directly.
$ws.Range("H3").Value = 531
$ws.Range("H4").Value = 84
$ws.Range("H5").Value = -50

# The "LAST UPDATE" cells hold plain text date labels (not real Excel
# dates) in the source workbook. Writing a date-shaped string straight to
# .Value would get auto-converted into a date serial, which would also
# pull in a brand-new number-format style. Instead, build the text via a
# text formula and then paste-special as values so the cell keeps its
# original (unaffected) style and ends up holding a literal string.
function Set-LiteralText {
    param($cell, [string]$text)
    $cell.Formula = '="' + $text + '"'
    $cell.Copy() | Out-Null
    $cell.PasteSpecial(-4163) | Out-Null  # xlPasteValues
}

Set-LiteralText $ws.Range("I3") "16-Sep-2025"
Set-LiteralText $ws.Range("I4") "16-Sep-2025"
Set-LiteralText $ws.Range("I5") "16-Sep-2025"

$excel.CutCopyMode = $false
